{"js": "// The author's edit adds a trailing space to the end of the\n// \"Example file: sgoop.py\" list item (the bullet right after the\n// \"Relevant parameters: sgoop.in_file, sgoop.rc_bin, sgoop.wells\" item,\n// under the \"Maximum Caliber\" heading).\n//\n// (The surrounding proofing-error/spell-check markers that Word had\n// inserted around hyphenated/dotted tokens like \"sgoop.in_file\" are also\n// cleaned up by Word automatically as part of any edit/re-save and are not\n// something an Office.js script can or needs to manipulate directly.)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Example file: sgoop.py\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Example file: sgoop.py\"');\n}\n\ntarget.insertText(\" \", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# The author's edit adds a trailing space to the end of the\n# \"Example file: sgoop.py\" list item (the bullet right after the\n# \"Relevant parameters: sgoop.in_file, sgoop.rc_bin, sgoop.wells\" item,\n# under the \"Maximum Caliber\" heading).\n#\n# (The surrounding proofing-error/spell-check markers that Word had\n# inserted around hyphenated/dotted tokens like \"sgoop.in_file\" are\n# regenerated/cleaned up by Word itself on edit/re-save and are not\n# something a COM automation script needs to manipulate directly.)\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Example file: sgoop.py\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph 'Example file: sgoop.py'\"\n}\n\n$target.Range.InsertAfter(\" \")\n"}
